# Author's intent (see commit message: "implementazione rule empty /
# adattamento al require_once di Data.php"): the "codice fiscale" (tax
# code) value that used to be written into E2 of Foglio1 is no longer
# produced by the data source, so the cell is cleared out entirely
# (its shared-string entry disappears from sharedStrings.xml as a
# result). The sheet's print setup was also touched (paper size /
# orientation now explicit instead of the "unset" 0 values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Remove the "codice fiscale" value that was in E2 - leave the cell
# genuinely empty rather than just blanking its text.
$ws.Range("E2").ClearContents()

# Update the page setup for the sheet (paper size / orientation are now
# explicitly set instead of the original placeholder zeros).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
